$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the address text in G6 (was the long "2375 Pennsylvania Av. NW, 20037 Washington DC")
$ws.Range("G6").Value = "Av. NW, 20037 Washington"

# Move the active cell selection to G6 (was G8)
$ws.Range("G6").Select()
